$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore window geometry (best effort - matches the committed view state)
$win = $wb.Windows.Item(1)
$win.Left = 6660
$win.Top = 3165
$win.Width = 21600
$win.Height = 11385

$ws.Range("A8").Value = "ΧΡΟΝΑΚΗΣ ΕΜΜΑΝΟΥΗΛ"
$ws.Range("E8").Value = 6937036009

$ws.Columns.Item(1).ColumnWidth = 22.67
$ws.Columns.Item(5).ColumnWidth = 8.33

$ws.Range("D9").Select()
